# Sample grade upload workbook: remove the sample student rows (test data),
# leaving just the header row and two blank, still-styled rows, and drop
# the mailto hyperlink that pointed at the removed sample e-mail addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Remove the mailto: hyperlink on A2 (and any others) before the cell
# contents referenced by it are cleared.
$ws.Hyperlinks.Delete()

# Wipe out the sample rows' data (A2:D4). This also drops row 4 entirely
# from the saved sheetData (it becomes fully empty) while leaving the
# cell-level styling on A2/A3 intact, and removes the now-unused shared
# strings for the sample names/e-mails.
$ws.Range("A2:D4").ClearContents()

# Match the author's last selection before saving.
$ws.Range("B2").Select()
